# Updated cure data (v42)
# Insert a new hospital row ("Roessingh, Centrum voor Revalidatie") into the
# alphabetically sorted list on Sheet1, keeping the existing sort order.
# The new row lands at row 45 (between "Rode Kruis Ziekenhuis" and
# "Slingeland Ziekenhuis (Santiz)"), pushing all following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 45, shifting rows 45..71 down to 46..72
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row with the new entry
$ws.Cells.Item(45, 1).Value = "Roessingh, Centrum voor Revalidatie"
$ws.Cells.Item(45, 2).Value = "voorlopig"

# Clear the selection / scroll position artifacts left over from editing so the
# sheet view resets to the top-left default, matching the saved state.
$ws.Range("A1").Select()
